# tweaks to systems & networks
$d = $word.ActiveDocument

# --- 1. Give the "Online resources for information security" bookmark a
#        human-readable name instead of the generated hash-style name. ---
$oldBookmarkName = "Xeba3c2f3420c339c3ff5bf9459ba3741edb2a9a"
$newBookmarkName = "online-resources-for-information-security"

$count = $d.Bookmarks.Count
$targetRange = $null
for ($i = 1; $i -le $count; $i++) {
    $bm = $d.Bookmarks.Item($i)
    if ($bm.Name -eq $oldBookmarkName) {
        $targetRange = $bm.Range
    }
}
if ($targetRange -ne $null) {
    $targetRange.Bookmarks.Add($newBookmarkName, $targetRange) | Out-Null
    # remove the old, auto-generated bookmark name now that the readable
    # one is in place
    for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
        $bm2 = $d.Bookmarks.Item($i)
        if ($bm2.Name -eq $oldBookmarkName) {
            $bm2.Delete()
        }
    }
}

# --- 2. The empty "Due" cells in the class-meetings table were missing
#        paragraph styling; give each of their (empty) paragraphs the
#        "Compact" style used throughout the rest of the table. ---
$table = $d.Tables.Item(1)
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    $cell = $row.Cells.Item(3)
    $cellRange = $cell.Range
    # an empty cell's Range.Text is just the end-of-cell marker (plus an
    # optional trailing CR) -- i.e. length <= 2 -- vs. real text otherwise
    if ($cellRange.Paragraphs.Count -eq 1 -and $cellRange.Text.Length -le 2) {
        $p = $cellRange.Paragraphs.Item(1)
        $p.Style = "Compact"
    }
}
